# Scheduled runner update: refresh computed market/profit figures on each
# class sheet (currentAveragePrice*, LevePrice*, LeveProfit* columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 75.25
$ws.Range("I5").Value = 83.666664
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 83.666664
$ws.Range("L5").Value = 50
$ws.Range("M5").Value = 31.333336
$ws.Range("N5").Value = -280
$ws.Range("H70").Value = 7289.8125
$ws.Range("J70").Value = 6992.0713
$ws.Range("L70").Value = 20976.2139
$ws.Range("N70").Value = -21516.2139
$ws.Range("H73").Value = 7289.8125
$ws.Range("J73").Value = 6992.0713
$ws.Range("L73").Value = 20976.2139
$ws.Range("N73").Value = -22848.2139
$ws.Range("H107").Value = 1454.3478
$ws.Range("I107").Value = 1650
$ws.Range("J107").Value = 1087.5
$ws.Range("K107").Value = 1650
$ws.Range("L107").Value = 1087.5
$ws.Range("M107").Value = 270
$ws.Range("N107").Value = -4927.5
$ws.Range("H126").Value = 39833.332
$ws.Range("J126").Value = 39833.332
$ws.Range("L126").Value = 39833.332
$ws.Range("N126").Value = -49713.332
$ws.Range("H137").Value = 10014562
$ws.Range("J137").Value = 3379.7
$ws.Range("L137").Value = 10139.1
$ws.Range("N137").Value = -15239.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7258.373
$ws.Range("I32").Value = 8000.9805
$ws.Range("K32").Value = 8000.9805
$ws.Range("M32").Value = -7713.9805
$ws.Range("H34").Value = 30009
$ws.Range("J34").Value = 30009
$ws.Range("L34").Value = 30009
$ws.Range("N34").Value = -30551
$ws.Range("H45").Value = 3108.75
$ws.Range("I45").Value = 2070
$ws.Range("K45").Value = 2070
$ws.Range("M45").Value = -1693
$ws.Range("H74").Value = 1316.1578
$ws.Range("I74").Value = 1306.2941
$ws.Range("J74").Value = 1400
$ws.Range("K74").Value = 1306.2941
$ws.Range("L74").Value = 1400
$ws.Range("M74").Value = -432.2941000000001
$ws.Range("N74").Value = -3148
$ws.Range("H77").Value = 1316.1578
$ws.Range("I77").Value = 1306.2941
$ws.Range("J77").Value = 1400
$ws.Range("K77").Value = 6531.4705
$ws.Range("L77").Value = 7000
$ws.Range("M77").Value = -2163.4705
$ws.Range("N77").Value = -15736
$ws.Range("H88").Value = 1741.3334
$ws.Range("J88").Value = 1699.1666
$ws.Range("L88").Value = 1699.1666
$ws.Range("N88").Value = -2511.1666
$ws.Range("H91").Value = 1741.3334
$ws.Range("J91").Value = 1699.1666
$ws.Range("L91").Value = 1699.1666
$ws.Range("N91").Value = -4507.1666
$ws.Range("H109").Value = 29833.334
$ws.Range("J109").Value = 29833.334
$ws.Range("L109").Value = 29833.334
$ws.Range("N109").Value = -32607.334
$ws.Range("H110").Value = 3338.2307
$ws.Range("I110").Value = 4388.75
$ws.Range("K110").Value = 4388.75
$ws.Range("M110").Value = -2343.75
$ws.Range("H122").Value = 5550
$ws.Range("I122").Value = 5253.3335
$ws.Range("K122").Value = 15760.0005
$ws.Range("M122").Value = -13310.0005
$ws.Range("H125").Value = 39777.777
$ws.Range("J125").Value = 39777.777
$ws.Range("L125").Value = 39777.777
$ws.Range("N125").Value = -49617.777
$ws.Range("H132").Value = 4808.8984
$ws.Range("I132").Value = 4029.8596
$ws.Range("J132").Value = 8509.333000000001
$ws.Range("K132").Value = 12089.5788
$ws.Range("L132").Value = 25527.999
$ws.Range("M132").Value = -9559.578799999999
$ws.Range("N132").Value = -30587.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 127402
$ws.Range("J86").Value = 2400
$ws.Range("L86").Value = 2400
$ws.Range("N86").Value = -4646
$ws.Range("H89").Value = 127402
$ws.Range("J89").Value = 2400
$ws.Range("L89").Value = 12000
$ws.Range("N89").Value = -23232
$ws.Range("H107").Value = 8271
$ws.Range("I107").Value = 8849.25
$ws.Range("K107").Value = 8849.25
$ws.Range("M107").Value = -6929.25
$ws.Range("H108").Value = 89742
$ws.Range("J108").Value = 89742
$ws.Range("L108").Value = 89742
$ws.Range("N108").Value = -97422
$ws.Range("H134").Value = 4048.6033
$ws.Range("I134").Value = 2715.0962
$ws.Range("K134").Value = 8145.2886
$ws.Range("M134").Value = -5610.2886

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 626250
$ws.Range("I4").Value = 1250000
$ws.Range("J4").Value = 2500
$ws.Range("K4").Value = 1250000
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = -1249888
$ws.Range("N4").Value = -2724
$ws.Range("H31").Value = 3124.4639
$ws.Range("I31").Value = 3201.26
$ws.Range("K31").Value = 3201.26
$ws.Range("M31").Value = -2906.26
$ws.Range("H34").Value = 3124.4639
$ws.Range("I34").Value = 3201.26
$ws.Range("K34").Value = 3201.26
$ws.Range("M34").Value = -2999.26
$ws.Range("H58").Value = 3158.8
$ws.Range("I58").Value = 2723.75
$ws.Range("K58").Value = 2723.75
$ws.Range("M58").Value = -2520.75
$ws.Range("H132").Value = 1182.3334
$ws.Range("I132").Value = 1182.3334
$ws.Range("K132").Value = 3547.0002
$ws.Range("M132").Value = -1017.0002
$ws.Range("H133").Value = 79975.336
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H136").Value = 3158.8
$ws.Range("I136").Value = 2723.75
$ws.Range("K136").Value = 8171.25
$ws.Range("M136").Value = -5621.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4993.6665
$ws.Range("J64").Value = 4993.6665
$ws.Range("L64").Value = 14980.9995
$ws.Range("N64").Value = -15520.9995
$ws.Range("H67").Value = 4993.6665
$ws.Range("J67").Value = 4993.6665
$ws.Range("L67").Value = 14980.9995
$ws.Range("N67").Value = -16852.9995
$ws.Range("H68").Value = 1537.6875
$ws.Range("J68").Value = 1857.5714
$ws.Range("L68").Value = 5572.7142
$ws.Range("N68").Value = -7194.7142
$ws.Range("H71").Value = 1537.6875
$ws.Range("J71").Value = 1857.5714
$ws.Range("L71").Value = 16718.1426
$ws.Range("N71").Value = -24830.1426
$ws.Range("H114").Value = 2441.4546
$ws.Range("J114").Value = 3595
$ws.Range("L114").Value = 10785
$ws.Range("N114").Value = -17293
$ws.Range("H122").Value = 3335.516
$ws.Range("I122").Value = 771.7273
$ws.Range("K122").Value = 6945.545700000001
$ws.Range("M122").Value = -4495.545700000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 4112251.5
$ws.Range("J24").Value = 29345.75
$ws.Range("L24").Value = 29345.75
$ws.Range("N24").Value = -29691.75
$ws.Range("H52").Value = 27554.334
$ws.Range("I52").Value = 21000
$ws.Range("J52").Value = 28373.625
$ws.Range("K52").Value = 21000
$ws.Range("L52").Value = 28373.625
$ws.Range("M52").Value = -20741
$ws.Range("N52").Value = -28891.625
$ws.Range("H57").Value = 15000
$ws.Range("I57").Value = 15000
$ws.Range("K57").Value = 15000
$ws.Range("M57").Value = -14180
$ws.Range("H80").Value = 113991
$ws.Range("I80").Value = 279977.75
$ws.Range("J80").Value = 3333.1667
$ws.Range("K80").Value = 279977.75
$ws.Range("L80").Value = 3333.1667
$ws.Range("M80").Value = -278979.75
$ws.Range("N80").Value = -5329.1667
$ws.Range("H83").Value = 113991
$ws.Range("I83").Value = 279977.75
$ws.Range("J83").Value = 3333.1667
$ws.Range("K83").Value = 1399888.75
$ws.Range("L83").Value = 16665.8335
$ws.Range("M83").Value = -1394896.75
$ws.Range("N83").Value = -26649.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2348.7585
$ws.Range("I93").Value = 1393.92
$ws.Range("K93").Value = 1393.92
$ws.Range("M93").Value = -145.9200000000001
$ws.Range("H136").Value = 2954.0667
$ws.Range("I136").Value = 2692
$ws.Range("K136").Value = 8076
$ws.Range("M136").Value = -5526
$ws.Range("H137").Value = 67571.36
$ws.Range("J137").Value = 69769.234
$ws.Range("L137").Value = 69769.234
$ws.Range("N137").Value = -79969.234

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 21776.666
$ws.Range("J94").Value = 21776.666
$ws.Range("L94").Value = 21776.666
$ws.Range("N94").Value = -23578.666
$ws.Range("H136").Value = 3162.8794
$ws.Range("I136").Value = 3223.9211
$ws.Range("J136").Value = 3046.9
$ws.Range("K136").Value = 9671.763300000001
$ws.Range("L136").Value = 9140.700000000001
$ws.Range("M136").Value = -7121.763300000001
$ws.Range("N136").Value = -14240.7
